# Manual Updates to reflect current status
#
# Applies:
#   - HandoutMaster / NotesMaster "datetimeFigureOut" date field text:
#       24/03/2022 -> 08/06/2022
#   - Slide 1, shape "Flowchart: Document 5":
#       body inset margins set to 0 (adds lIns="0" rIns="0")
#       run text "Defines Target Application" ->
#                "Target (Device, Board) Build (Debug, Release) "
#   - Slide 1, shape "Flowchart: Document 37":
#       run text "csettings.yml" -> "cdefault.yml"
#       run text "(Toolchain, Device)" -> "(Toolchain)"

$p = $ppt.ActivePresentation

function Replace-Substring($textRange, [string]$oldStr, [string]$newStr) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldStr)
    if ($idx -lt 0) {
        throw "substring not found: [$oldStr]"
    }
    $sub = $textRange.Characters($idx + 1, $oldStr.Length)
    $sub.Text = $newStr
}

# --- Handout master / notes master date placeholders ---
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Text = "08/06/2022"

$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "08/06/2022"

# --- Slide 1 ---
$s = $p.Slides.Item(1)

# "Flowchart: Document 5" - Target Application box
$shTarget = $s.Shapes.Item("Flowchart: Document 5")
$shTarget.TextFrame.MarginLeft = 0
$shTarget.TextFrame.MarginRight = 0
Replace-Substring $shTarget.TextFrame.TextRange "Defines Target Application" "Target (Device, Board) Build (Debug, Release) "

# "Flowchart: Document 37" - csettings.yml box
$shSettings = $s.Shapes.Item("Flowchart: Document 37")
Replace-Substring $shSettings.TextFrame.TextRange "csettings.yml" "cdefault.yml"
Replace-Substring $shSettings.TextFrame.TextRange "(Toolchain, Device)" "(Toolchain)"
